# Apply "F" column (想去人数) updates across sheets "展览" and "演出",
# plus the mirrored rows on the aggregate "全部类型" sheet.
# Values below are (row -> new value) per sheet, matching the commit's
# regenerated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

function Set-FValues {
    # NOTE: use positional params here -- passing Unicode sheet names via
    # named parameters (-SheetName "...") silently truncates them to an
    # empty string in this runtime, which then makes Worksheets.Item throw
    # "subscript out of range". Positional binding works correctly.
    param($SheetName, $RowValues)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Range("F$row").Value = $RowValues[$row]
    }
}

# Sheet "展览" (sheet1)
Set-FValues "展览" @{
    4  = 5907
    5  = 5907
    7  = 2958
    8  = 1271
    13 = 216
    14 = 4278
    15 = 4278
    18 = 107
    22 = 6474
    23 = 6474
    24 = 227
    26 = 292
    28 = 1225
    31 = 1624
    33 = 1856
    34 = 5951
    35 = 105
    36 = 19
    39 = 397
    40 = 4071
    42 = 185
    43 = 80
    45 = 2400
    46 = 21
    47 = 44
    50 = 311
    52 = 19
}

# Sheet "演出" (sheet2)
Set-FValues "演出" @{
    5 = 99
}

# Sheet "全部类型" (sheet4) - aggregate of all other sheets
Set-FValues "全部类型" @{
    7  = 2958
    8  = 1271
    13 = 216
    14 = 4278
    15 = 4278
    18 = 107
    22 = 6474
    23 = 6474
    24 = 227
    27 = 1225
    28 = 99
    30 = 1624
    33 = 1856
    35 = 5951
    36 = 105
    37 = 19
    40 = 397
    41 = 4071
    42 = 185
    43 = 80
    47 = 2400
    48 = 21
    49 = 44
    51 = 311
}
